# Auto-generated COM script to apply the 2025-12-01 18:30 JST scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# 1) Drop the old hyperlink relationships that live on column F (rows 2-14);
#    the cell text/style is rewritten below and fresh hyperlinks are added afterwards.
$ws.Range("F2:F14").Hyperlinks.Delete()

# 2) Rewrite every data row (2-23) with the refreshed scrape contents.
# row 2
$ws.Cells.Item(2,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(2,2).Value = "【急募】生成AI×業務効率化の実装を支援するエンジニア募集"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5444662"
$ws.Cells.Item(2,7).Value = 385
$ws.Cells.Item(2,8).Value = "🔥AI,Ai ◆効率化"

# row 3
$ws.Cells.Item(3,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(3,2).Value = "B2B向け生成AIサービス(チャット・RAG)の新規開発プロジェクト推進を支援してくださるPM募集"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,6).Value = "https://www.lancers.jp/work/detail/5445154"
$ws.Cells.Item(3,7).Value = 368
$ws.Cells.Item(3,8).Value = "🔥AI,Ai ◆開発"

# row 4
$ws.Cells.Item(4,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(4,2).Value = "建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集"
$ws.Cells.Item(4,3).Value = "システム開発"
$ws.Cells.Item(4,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(4,5).Value = "期限情報なし"
$ws.Cells.Item(4,6).Value = "https://www.lancers.jp/work/detail/5434128"
$ws.Cells.Item(4,7).Value = 368
$ws.Cells.Item(4,8).Value = "🔥AI,Ai ◆開発"

# row 5
$ws.Cells.Item(5,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(5,2).Value = "法人向け生成AIサービス(RAG・議事録機能)の設計・開発を支援エンジニア募集(AI/バックエンド)"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5445159"
$ws.Cells.Item(5,7).Value = 368
$ws.Cells.Item(5,8).Value = "🔥AI,Ai ◆開発"

# row 6
$ws.Cells.Item(6,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(6,2).Value = "企業のMicrosoft Copilot導入・活用支援AIコンサルタント募集(研修講師・メンター)"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5434363"
$ws.Cells.Item(6,7).Value = 348
$ws.Cells.Item(6,8).Value = "🔥AI,Ai ◆コンサル"

# row 7
$ws.Cells.Item(7,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(7,2).Value = "【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像認識/動画解析)"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,6).Value = "https://www.lancers.jp/work/detail/5439158"
$ws.Cells.Item(7,7).Value = 303
$ws.Cells.Item(7,8).Value = "🔥AI,Ai"

# row 8
$ws.Cells.Item(8,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(8,2).Value = "製造業のR&D支援!「プロセスデータ解析」「音響異常検知」のAIエンジニア募集"
$ws.Cells.Item(8,3).Value = "システム開発"
$ws.Cells.Item(8,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(8,5).Value = "期限情報なし"
$ws.Cells.Item(8,6).Value = "https://www.lancers.jp/work/detail/5439165"
$ws.Cells.Item(8,7).Value = 303
$ws.Cells.Item(8,8).Value = "🔥AI,Ai"

# row 9
$ws.Cells.Item(9,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(9,2).Value = "外部WEB予約サイト一元管理システム開発|長期保守パートナー募集"
$ws.Cells.Item(9,3).Value = "システム開発"
$ws.Cells.Item(9,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(9,5).Value = "期限情報なし"
$ws.Cells.Item(9,6).Value = "https://www.lancers.jp/work/detail/5444378"
$ws.Cells.Item(9,7).Value = 170
$ws.Cells.Item(9,8).Value = "◆開発,システム開発 ◇サイト"

# row 10
$ws.Cells.Item(10,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(10,2).Value = "【日本人限定/継続案件】Node.jsエンジニア募集(スクレイピング機能開発)"
$ws.Cells.Item(10,3).Value = "システム開発"
$ws.Cells.Item(10,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(10,5).Value = "期限情報なし"
$ws.Cells.Item(10,6).Value = "https://www.lancers.jp/work/detail/5444489"
$ws.Cells.Item(10,7).Value = 155
$ws.Cells.Item(10,8).Value = "◆開発,Node.js"

# row 11
$ws.Cells.Item(11,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(11,2).Value = "【効率化】Microsoft 365で英国イベントの出展者登録や情報管理を簡素化"
$ws.Cells.Item(11,3).Value = "システム開発"
$ws.Cells.Item(11,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(11,5).Value = "期限情報なし"
$ws.Cells.Item(11,6).Value = "https://www.lancers.jp/work/detail/5445148"
$ws.Cells.Item(11,7).Value = 103
$ws.Cells.Item(11,8).Value = "◆効率化 ◇管理"

# row 12
$ws.Cells.Item(12,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(12,2).Value = "【Excelでのマクロ作成】リサーチツールの作成【スクレイピング】"
$ws.Cells.Item(12,3).Value = "システム開発"
$ws.Cells.Item(12,4).Value = "1,000 ~ 5,000 円 / 固定"
$ws.Cells.Item(12,5).Value = "期限情報なし"
$ws.Cells.Item(12,6).Value = "https://www.lancers.jp/work/detail/5445173"
$ws.Cells.Item(12,7).Value = 100
$ws.Cells.Item(12,8).Value = "◆ツール,スクレイピング"

# row 13
$ws.Cells.Item(13,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(13,2).Value = "【Excelでのマクロ作成】リサーチツールの作成【スクレイピング】"
$ws.Cells.Item(13,3).Value = "システム開発"
$ws.Cells.Item(13,4).Value = "1,000 ~ 5,000 円 / 固定"
$ws.Cells.Item(13,5).Value = "期限情報なし"
$ws.Cells.Item(13,6).Value = "https://www.lancers.jp/work/detail/5445149"
$ws.Cells.Item(13,7).Value = 100
$ws.Cells.Item(13,8).Value = "◆ツール,スクレイピング"

# row 14
$ws.Cells.Item(14,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(14,2).Value = "【急募】Googledriveのロール管理・共有設定の専門家募集"
$ws.Cells.Item(14,3).Value = "システム開発"
$ws.Cells.Item(14,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(14,5).Value = "期限情報なし"
$ws.Cells.Item(14,6).Value = "https://www.lancers.jp/work/detail/5444395"
$ws.Cells.Item(14,7).Value = 38
$ws.Cells.Item(14,8).Value = "◇管理"

# row 15
$ws.Cells.Item(15,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(15,2).Value = "【急募】freee勤怠管理plus用シフトCSV作成依頼"
$ws.Cells.Item(15,3).Value = "システム開発"
$ws.Cells.Item(15,4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(15,5).Value = "期限情報なし"
$ws.Cells.Item(15,6).Value = "https://www.lancers.jp/work/detail/5445210"
$ws.Cells.Item(15,7).Value = 30
$ws.Cells.Item(15,8).Value = "◇管理"

# row 16
$ws.Cells.Item(16,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(16,2).Value = "進行管理およびチームディレクションを担当"
$ws.Cells.Item(16,3).Value = "システム開発"
$ws.Cells.Item(16,4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(16,5).Value = "期限情報なし"
$ws.Cells.Item(16,6).Value = "https://www.lancers.jp/work/detail/5418064"
$ws.Cells.Item(16,7).Value = 30
$ws.Cells.Item(16,8).Value = "◇管理"

# row 17
$ws.Cells.Item(17,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(17,2).Value = "【急募】ネットワーク環境の確認と最適化をお願いします"
$ws.Cells.Item(17,3).Value = "システム開発"
$ws.Cells.Item(17,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(17,5).Value = "期限情報なし"
$ws.Cells.Item(17,6).Value = "https://www.lancers.jp/work/detail/5445215"
$ws.Cells.Item(17,7).Value = 18
$ws.Cells.Item(17,8).ClearContents()

# row 18
$ws.Cells.Item(18,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(18,2).Value = "【USPTOへ特許申請】dAppsのwebsite制作・Velodromeでのプレセール知識が必須"
$ws.Cells.Item(18,3).Value = "システム開発"
$ws.Cells.Item(18,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(18,5).Value = "期限情報なし"
$ws.Cells.Item(18,6).Value = "https://www.lancers.jp/work/detail/5445167"
$ws.Cells.Item(18,7).Value = 18
$ws.Cells.Item(18,8).ClearContents()

# row 19
$ws.Cells.Item(19,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(19,2).Value = "DAO構築。ブロックチェーンとスマートコントラクトの専門家募集"
$ws.Cells.Item(19,3).Value = "システム開発"
$ws.Cells.Item(19,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(19,5).Value = "期限情報なし"
$ws.Cells.Item(19,6).Value = "https://www.lancers.jp/work/detail/5445105"
$ws.Cells.Item(19,7).Value = 18
$ws.Cells.Item(19,8).ClearContents()

# row 20
$ws.Cells.Item(20,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(20,2).Value = "ホスティング業務を担当してくれる方を探しています!"
$ws.Cells.Item(20,3).Value = "システム開発"
$ws.Cells.Item(20,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(20,5).Value = "期限情報なし"
$ws.Cells.Item(20,6).Value = "https://www.lancers.jp/work/detail/5445080"
$ws.Cells.Item(20,7).Value = 18
$ws.Cells.Item(20,8).ClearContents()

# row 21
$ws.Cells.Item(21,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(21,2).Value = "【急募】Amazonフラットファイル(ブラウズノード検証)に詳しい方を探しています"
$ws.Cells.Item(21,3).Value = "システム開発"
$ws.Cells.Item(21,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(21,5).Value = "期限情報なし"
$ws.Cells.Item(21,6).Value = "https://www.lancers.jp/work/detail/5444446"
$ws.Cells.Item(21,7).Value = 18
$ws.Cells.Item(21,8).ClearContents()

# row 22
$ws.Cells.Item(22,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(22,2).Value = "【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え"
$ws.Cells.Item(22,3).Value = "システム開発"
$ws.Cells.Item(22,4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(22,5).Value = "期限情報なし"
$ws.Cells.Item(22,6).Value = "https://www.lancers.jp/work/detail/5443568"
$ws.Cells.Item(22,7).Value = 13
$ws.Cells.Item(22,8).ClearContents()

# row 23
$ws.Cells.Item(23,1).Value = "2025-12-01 18:30:10"
$ws.Cells.Item(23,2).Value = "comfyui(paperspace)でエロ動画のループ物を作成したいです。その方法を教えてください"
$ws.Cells.Item(23,3).Value = "システム開発"
$ws.Cells.Item(23,4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(23,5).Value = "期限情報なし"
$ws.Cells.Item(23,6).Value = "https://www.lancers.jp/work/detail/5444370"
$ws.Cells.Item(23,7).Value = 10
$ws.Cells.Item(23,8).ClearContents()

# 3) Re-create the hyperlinks on column F for every data row, then restore the
#    shared "Hyperlink" cell style (Hyperlinks.Add otherwise registers a spare style).
$ws.Hyperlinks.Add($ws.Cells.Item(2,6), "https://www.lancers.jp/work/detail/5444662") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3,6), "https://www.lancers.jp/work/detail/5445154") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4,6), "https://www.lancers.jp/work/detail/5434128") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5,6), "https://www.lancers.jp/work/detail/5445159") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(6,6), "https://www.lancers.jp/work/detail/5434363") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(7,6), "https://www.lancers.jp/work/detail/5439158") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(8,6), "https://www.lancers.jp/work/detail/5439165") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(9,6), "https://www.lancers.jp/work/detail/5444378") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(10,6), "https://www.lancers.jp/work/detail/5444489") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(11,6), "https://www.lancers.jp/work/detail/5445148") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(12,6), "https://www.lancers.jp/work/detail/5445173") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(13,6), "https://www.lancers.jp/work/detail/5445149") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(14,6), "https://www.lancers.jp/work/detail/5444395") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(15,6), "https://www.lancers.jp/work/detail/5445210") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(16,6), "https://www.lancers.jp/work/detail/5418064") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(17,6), "https://www.lancers.jp/work/detail/5445215") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(18,6), "https://www.lancers.jp/work/detail/5445167") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(19,6), "https://www.lancers.jp/work/detail/5445105") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(20,6), "https://www.lancers.jp/work/detail/5445080") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(21,6), "https://www.lancers.jp/work/detail/5444446") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(22,6), "https://www.lancers.jp/work/detail/5443568") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(23,6), "https://www.lancers.jp/work/detail/5444370") | Out-Null
$ws.Cells.Item(2,6).Style = "Hyperlink"
$ws.Cells.Item(3,6).Style = "Hyperlink"
$ws.Cells.Item(4,6).Style = "Hyperlink"
$ws.Cells.Item(5,6).Style = "Hyperlink"
$ws.Cells.Item(6,6).Style = "Hyperlink"
$ws.Cells.Item(7,6).Style = "Hyperlink"
$ws.Cells.Item(8,6).Style = "Hyperlink"
$ws.Cells.Item(9,6).Style = "Hyperlink"
$ws.Cells.Item(10,6).Style = "Hyperlink"
$ws.Cells.Item(11,6).Style = "Hyperlink"
$ws.Cells.Item(12,6).Style = "Hyperlink"
$ws.Cells.Item(13,6).Style = "Hyperlink"
$ws.Cells.Item(14,6).Style = "Hyperlink"
$ws.Cells.Item(15,6).Style = "Hyperlink"
$ws.Cells.Item(16,6).Style = "Hyperlink"
$ws.Cells.Item(17,6).Style = "Hyperlink"
$ws.Cells.Item(18,6).Style = "Hyperlink"
$ws.Cells.Item(19,6).Style = "Hyperlink"
$ws.Cells.Item(20,6).Style = "Hyperlink"
$ws.Cells.Item(21,6).Style = "Hyperlink"
$ws.Cells.Item(22,6).Style = "Hyperlink"
$ws.Cells.Item(23,6).Style = "Hyperlink"

Write-Host "applied 2025-12-01 18:30:10 refresh"
